$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")

# Test data rows 10 and 11: the PCH800 / PLX800 columns (L, M) were being
# shown as the "NA" text value - they are filtered/boolean columns, so
# they should hold the boolean FALSE instead.
$ws.Range("L10").Value = $false
$ws.Range("M10").Value = $false
$ws.Range("L11").Value = $false
$ws.Range("M11").Value = $false

# Reflect the author's updated view/selection on the "Add Devices" sheet.
$ws.Activate()
$ws.Range("A10").Select()
